# Update database: drop the oldest period (1396/12) and its publish date,
# shift every period one column to the left, and bring in the newly
# published 1401/12 period with its own publish dates.
# Also fixes the D15 placeholder ("-") to a real numeric 0, matching the
# new read_price algorithm that always emits a number for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers -------------------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ----------------------------------------------------
$ws.Range("D9").Value = "1399-03-21 (10)"
$ws.Range("E9").Value = "1400-02-31 (10)"
$ws.Range("F9").Value = "1401-02-31 (11)"
$ws.Range("G9").Value = "1402-02-25 (10)"
$ws.Range("H9").Value = "1402-02-25 (2)"

# --- Row 11: فروش ------------------------------------------------------------
$ws.Range("D11").Value = 2812032
$ws.Range("E11").Value = 5177955
$ws.Range("F11").Value = 8083748
$ws.Range("G11").Value = 12824117
$ws.Range("H11").Value = 19977759

# --- Row 12: بهای تمام شده کالای فروش رفته -----------------------------------
$ws.Range("D12").Value = -1669621
$ws.Range("E12").Value = -3085138
$ws.Range("F12").Value = -3789830
$ws.Range("G12").Value = -6208752
$ws.Range("H12").Value = -9443386

# --- Row 13: سود (زیان) ناخالص ------------------------------------------------
$ws.Range("D13").Value = 1142411
$ws.Range("E13").Value = 2092817
$ws.Range("F13").Value = 4293918
$ws.Range("G13").Value = 6615365
$ws.Range("H13").Value = 10534373

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی -------------------------------
$ws.Range("D14").Value = -273056
$ws.Range("E14").Value = -969095
$ws.Range("F14").Value = -1030517
$ws.Range("G14").Value = -1779303
$ws.Range("H14").Value = -2619163

# --- Row 15: هزینه کاهش ارزش دریافتنی‌‏ها (هزینه استثنایی) -------------------
# D15 used to hold the placeholder text "-"; the new read_price algorithm
# now always resolves it to a real number.
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---------------------------
$ws.Range("D16").Value = 41357
$ws.Range("E16").Value = 73369
$ws.Range("F16").Value = 289564
$ws.Range("G16").Value = 165861
$ws.Range("H16").Value = 1736051

# --- Row 17: سود (زیان) عملیاتی ----------------------------------------------
$ws.Range("D17").Value = 910712
$ws.Range("E17").Value = 1197091
$ws.Range("F17").Value = 3552965
$ws.Range("G17").Value = 5001923
$ws.Range("H17").Value = 9651261

# --- Row 18: هزینه های مالی ---------------------------------------------------
$ws.Range("D18").Value = -229102
$ws.Range("E18").Value = -134733
$ws.Range("F18").Value = -52717
$ws.Range("G18").Value = -4617
$ws.Range("H18").Value = -18198

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی -------------------------
$ws.Range("D19").Value = 46374
$ws.Range("E19").Value = 24776
$ws.Range("F19").Value = 357261
$ws.Range("G19").Value = 343956
$ws.Range("H19").Value = 843642

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---------------
$ws.Range("D20").Value = 727984
$ws.Range("E20").Value = 1087134
$ws.Range("F20").Value = 3857509
$ws.Range("G20").Value = 5341262
$ws.Range("H20").Value = 10476705

# --- Row 21: مالیات ------------------------------------------------------------
$ws.Range("D21").Value = -33378
$ws.Range("E21").Value = -58672
$ws.Range("F21").Value = -174187
$ws.Range("G21").Value = -111473
$ws.Range("H21").Value = -411036

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ------------------------------
$ws.Range("D22").Value = 694606
$ws.Range("E22").Value = 1028462
$ws.Range("F22").Value = 3683322
$ws.Range("G22").Value = 5229789
$ws.Range("H22").Value = 10065669

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی --------------------
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# --- Row 24: سود (زیان) خالص ---------------------------------------------------
$ws.Range("D24").Value = 694606
$ws.Range("E24").Value = 1028462
$ws.Range("F24").Value = 3683322
$ws.Range("G24").Value = 5229789
$ws.Range("H24").Value = 10065669

# --- Row 25: سود هر سهم پس از کسر مالیات ---------------------------------------
$ws.Range("D25").Value = 1069
$ws.Range("E25").Value = 735
$ws.Range("F25").Value = 2631
$ws.Range("G25").Value = 3736
$ws.Range("H25").Value = 7190

# --- Row 26: سرمایه -------------------------------------------------------------
$ws.Range("D26").Value = 650000
$ws.Range("E26").Value = 1400000
$ws.Range("F26").Value = 1400000
$ws.Range("G26").Value = 1400000
$ws.Range("H26").Value = 1400000

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه ------------------------------------
$ws.Range("D27").Value = 496
$ws.Range("E27").Value = 735
$ws.Range("F27").Value = 2631
$ws.Range("G27").Value = 3736
$ws.Range("H27").Value = 7190
